$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column H ("property_category") before the current date column,
# pushing date / legislator_name / legislator_id one column to the right.
$ws.Columns.Item(8).Insert()

# Header for the new column.
$ws.Cells.Item(1, 8).Value = "property_category"

# Fill the new column with "stock" for every existing data row.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
